# Build OWL file for Behaviour module
#
# 1. Append "; function [BFO:0000034]" to the BFO class-list cell (D5 on the
#    "Imports" sheet) that enumerates the BFO classes used as roots/parents.
# 2. Move the active selection on the sheet from D6 to E5 (this is what the
#    author's click/selection ended up at when the workbook was saved).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cell = $ws.Range("D5")
$cell.Value = $cell.Text + "; function [BFO:0000034]"

$ws.Range("E5").Select()
